# Update on 20 Jan
# Append the remaining column descriptions (GarageCars .. SaleCondition) to the
# "Columns_Analysis" worksheet, continuing the existing table (rows 59-77).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(59, 'GarageCars', 'Numerical', 'Identified by Numbers', 'Size of garage in car capacity'),
    @(60, 'GarageArea', 'Numerical', 'Identified by Numbers', 'Size of garage in square feet'),
    @(61, 'GarageQual', 'Categorical', 'Identified by Strings', 'Garage quality'),
    @(62, 'GarageCond', 'Categorical', 'Identified by Strings', 'Garage condition'),
    @(63, 'PavedDrive', 'Categorical', 'Identified by Strings', 'Paved driveway'),
    @(64, 'WoodDeckSF', 'Numerical', 'Identified by Numbers', 'Wood deck area in square feet'),
    @(65, 'OpenPorchSF', 'Numerical', 'Identified by Numbers', 'Open porch area in square feet'),
    @(66, 'EnclosedPorch', 'Numerical', 'Identified by Numbers', 'Enclosed porch area in square feet'),
    @(67, '3SsnPorch', 'Numerical', 'Identified by Numbers', 'Three season porch area in square feet'),
    @(68, 'ScreenPorch', 'Numerical', 'Identified by Numbers', 'Screen porch area in square feet'),
    @(69, 'PoolArea', 'Numerical', 'Identified by Numbers', 'Pool area in square feet'),
    @(70, 'PoolQC', 'Categorical', 'Identified by Strings', 'Pool quality'),
    @(71, 'Fence', 'Categorical', 'Identified by Strings', 'Fence quality'),
    @(72, 'MiscFeature', 'Categorical', 'Identified by Strings', 'Miscellaneous feature not covered in other categories'),
    @(73, 'MiscVal', 'Numerical', 'Identified by Numbers', '$Value of miscellaneous feature'),
    @(74, 'MoSold', 'Numerical', 'Identified by Numbers', 'Month Sold (MM)'),
    @(75, 'YrSold', 'Numerical', 'Identified by Numbers', 'Year Sold (YYYY)'),
    @(76, 'SaleType', 'Categorical', 'Identified by Strings', 'Type of sale'),
    @(77, 'SaleCondition', 'Categorical', 'Identified by Strings', 'Condition of sale')
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
}

# Scroll the view down to the new data and move the active selection,
# matching where the author left off editing.
$excel.ActiveWindow.ScrollRow = 59
$ws.Range("D80").Select()
